$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values (e.g. "245.87")
# are not auto-converted to numbers by Excel COM type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.240.38"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.876.44"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "245.87"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").Value = "0.677"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").Value = "43.52"
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("D9").Value = "0.357"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "0.0739"
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("D12").Value = "0.0974"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "13.45"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").Value = "2.148.78"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "0.768"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.91"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.866.18"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "35.243.54"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "73.55"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "243.89"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "12.81"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "5.01"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "2.66"
$ws.Range("E24").Value = "  +9.80%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "2.18"
$ws.Range("E26").Value = "  -5.24%  "
$ws.Range("D27").Value = "165.36"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "8.56"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "18.25"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").Value = "4.28"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "0.0589"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.16"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  -11.78%  "
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").Value = "  -9.46%  "
$ws.Range("D37").Value = "0.846"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.0736"
$ws.Range("E38").Value = "  +9.40%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "0.0217"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "96.18"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").Value = "1.305.32"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").Value = "0.0799"
$ws.Range("E46").Value = "  +5.42%  "
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "2.72"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").Value = "11.73"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").Value = "6.26"
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "41.97"
$ws.Range("E51").Value = "  -2.34%  "

# Restore default (Normal) style on column D so only the number format used
# for text-coercion is cleared, matching original unstyled data cells.
$ws.Range("D2:D51").Style = "Normal"
